$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 1438905.8
$ws.Range("I33").Value = 2874262.2
$ws.Range("K33").Value = 2874262.2
$ws.Range("M33").Value = -2874033.2
# Row 132
$ws.Range("H132").Value = 1771.2963
$ws.Range("I132").Value = 1647.1154
$ws.Range("K132").Value = 4941.3462
$ws.Range("M132").Value = -2411.3462
# Row 137
$ws.Range("H137").Value = 4153.4585
$ws.Range("I137").Value = 1933.6786
$ws.Range("J137").Value = 7261.15
$ws.Range("K137").Value = 5801.0358
$ws.Range("L137").Value = 21783.45
$ws.Range("M137").Value = -3251.0358
$ws.Range("N137").Value = -26883.45
# Row 138
$ws.Range("H138").Value = 2181.8809
$ws.Range("I138").Value = 1391.28
$ws.Range("J138").Value = 3344.5293
$ws.Range("K138").Value = 4173.84
$ws.Range("L138").Value = 10033.5879
$ws.Range("M138").Value = 966.1599999999999
$ws.Range("N138").Value = -20313.5879

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3645153.5
$ws.Range("I32").Value = 719117.2
$ws.Range("J32").Value = 18275336
$ws.Range("K32").Value = 719117.2
$ws.Range("L32").Value = 18275336
$ws.Range("M32").Value = -718830.2
$ws.Range("N32").Value = -18275910
# Row 45
$ws.Range("H45").Value = 12097.207
$ws.Range("I45").Value = 9721.583000000001
$ws.Range("J45").Value = 23500.2
$ws.Range("K45").Value = 9721.583000000001
$ws.Range("L45").Value = 23500.2
$ws.Range("M45").Value = -9344.583000000001
$ws.Range("N45").Value = -24254.2
# Row 97
$ws.Range("H97").Value = 36719.035
$ws.Range("I97").Value = 8692.615
$ws.Range("J97").Value = 218890.75
$ws.Range("K97").Value = 8692.615
$ws.Range("L97").Value = 218890.75
$ws.Range("M97").Value = -8196.615
$ws.Range("N97").Value = -219882.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4970.278
$ws.Range("I20").Value = 4076.1428
$ws.Range("J20").Value = 8099.75
$ws.Range("K20").Value = 4076.1428
$ws.Range("L20").Value = 8099.75
$ws.Range("M20").Value = -3829.1428
$ws.Range("N20").Value = -8593.75
# Row 94
$ws.Range("H94").Value = 86852.336
$ws.Range("I94").Value = 2420
$ws.Range("J94").Value = 192392.75
$ws.Range("K94").Value = 2420
$ws.Range("L94").Value = 192392.75
$ws.Range("M94").Value = -1969
$ws.Range("N94").Value = -193294.75
# Row 105
$ws.Range("H105").Value = 17154.807
$ws.Range("I105").Value = 4227.909
$ws.Range("J105").Value = 48753.89
$ws.Range("K105").Value = 4227.909
$ws.Range("L105").Value = 48753.89
$ws.Range("M105").Value = -2480.909
$ws.Range("N105").Value = -52247.89
# Row 107
$ws.Range("H107").Value = 2345.3333
$ws.Range("I107").Value = 2321.2334
$ws.Range("J107").Value = 2465.8333
$ws.Range("K107").Value = 2321.2334
$ws.Range("L107").Value = 2465.8333
$ws.Range("M107").Value = -401.2334000000001
$ws.Range("N107").Value = -6305.8333
# Row 134
$ws.Range("H134").Value = 2948.0227
$ws.Range("I134").Value = 2892.475
$ws.Range("K134").Value = 8677.424999999999
$ws.Range("M134").Value = -6142.424999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2978
$ws.Range("I31").Value = 1208.762
$ws.Range("J31").Value = 7622.25
$ws.Range("K31").Value = 1208.762
$ws.Range("L31").Value = 7622.25
$ws.Range("M31").Value = -913.7619999999999
$ws.Range("N31").Value = -8212.25
# Row 34
$ws.Range("H34").Value = 2978
$ws.Range("I34").Value = 1208.762
$ws.Range("J34").Value = 7622.25
$ws.Range("K34").Value = 1208.762
$ws.Range("L34").Value = 7622.25
$ws.Range("M34").Value = -1006.762
$ws.Range("N34").Value = -8026.25
# Row 132
$ws.Range("H132").Value = 2387.4
$ws.Range("I132").Value = 2298.2334
$ws.Range("K132").Value = 6894.7002
$ws.Range("M132").Value = -4364.7002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1314.4
$ws.Range("I68").Value = 1358.1428
$ws.Range("J68").Value = 1290.8462
$ws.Range("K68").Value = 4074.4284
$ws.Range("L68").Value = 3872.5386
$ws.Range("M68").Value = -3263.4284
$ws.Range("N68").Value = -5494.5386
# Row 71
$ws.Range("H71").Value = 1314.4
$ws.Range("I71").Value = 1358.1428
$ws.Range("J71").Value = 1290.8462
$ws.Range("K71").Value = 12223.2852
$ws.Range("L71").Value = 11617.6158
$ws.Range("M71").Value = -8167.2852
$ws.Range("N71").Value = -19729.6158
# Row 76
$ws.Range("H76").Value = 3999
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 3999
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
# Row 81
$ws.Range("H81").Value = 45550104
$ws.Range("J81").Value = 125009500
$ws.Range("L81").Value = 375028500
$ws.Range("N81").Value = -375030746
# Row 84
$ws.Range("H84").Value = 45550104
$ws.Range("J84").Value = 125009500
$ws.Range("L84").Value = 1125085500
$ws.Range("N84").Value = -1125096732
# Row 87
$ws.Range("H87").Value = 15006.5
$ws.Range("J87").Value = 29999
$ws.Range("L87").Value = 89997
$ws.Range("N87").Value = -92493
# Row 90
$ws.Range("H90").Value = 15006.5
$ws.Range("J90").Value = 29999
$ws.Range("L90").Value = 269991
$ws.Range("N90").Value = -282471
# Row 93
$ws.Range("H93").Value = 128011
$ws.Range("J93").Value = 3998.1667
$ws.Range("L93").Value = 11994.5001
$ws.Range("N93").Value = -15738.5001
# Row 107
$ws.Range("H107").Value = 706.5599999999999
$ws.Range("I107").Value = 688.2
$ws.Range("J107").Value = 711.15
$ws.Range("K107").Value = 2064.6
$ws.Range("L107").Value = 2133.45
$ws.Range("M107").Value = -144.6000000000004
$ws.Range("N107").Value = -5973.45

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 22739566
$ws.Range("I107").Value = 998.9231
$ws.Range("K107").Value = 998.9231
$ws.Range("M107").Value = 921.0769
# Row 132
$ws.Range("H132").Value = 4904.552
$ws.Range("I132").Value = 5218.4585
$ws.Range("J132").Value = 3397.8
$ws.Range("K132").Value = 15655.3755
$ws.Range("L132").Value = 10193.4
$ws.Range("M132").Value = -13125.3755
$ws.Range("N132").Value = -15253.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1122.1177
$ws.Range("I16").Value = 1058.4
$ws.Range("J16").Value = 1600
$ws.Range("K16").Value = 1058.4
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = -888.4000000000001
$ws.Range("N16").Value = -1940
# Row 93
$ws.Range("H93").Value = 56997.668
$ws.Range("I93").Value = 1717.6666
$ws.Range("J93").Value = 112277.664
$ws.Range("K93").Value = 1717.6666
$ws.Range("L93").Value = 112277.664
$ws.Range("M93").Value = -469.6666
$ws.Range("N93").Value = -114773.664
# Row 132
$ws.Range("H132").Value = 5989.533
$ws.Range("I132").Value = 6363.222
$ws.Range("J132").Value = 5429
$ws.Range("K132").Value = 19089.666
$ws.Range("L132").Value = 16287
$ws.Range("M132").Value = -16559.666
$ws.Range("N132").Value = -21347

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 8735.75
$ws.Range("I81").Value = 8221.75
$ws.Range("J81").Value = 9249.75
$ws.Range("K81").Value = 16443.5
$ws.Range("L81").Value = 18499.5
$ws.Range("M81").Value = -15382.5
$ws.Range("N81").Value = -20621.5
# Row 84
$ws.Range("H84").Value = 8735.75
$ws.Range("I84").Value = 8221.75
$ws.Range("J84").Value = 9249.75
$ws.Range("K84").Value = 82217.5
$ws.Range("L84").Value = 92497.5
$ws.Range("M84").Value = -76913.5
$ws.Range("N84").Value = -103105.5
# Row 107
$ws.Range("H107").Value = 62572920
$ws.Range("I107").Value = 2615.8333
$ws.Range("J107").Value = 250283820
$ws.Range("K107").Value = 7847.499899999999
$ws.Range("L107").Value = 750851460
$ws.Range("M107").Value = -5927.499899999999
$ws.Range("N107").Value = -750855300
# Row 132
$ws.Range("H132").Value = 1072243.6
$ws.Range("I132").Value = 1279062.5
$ws.Range("J132").Value = 3679.3333
$ws.Range("K132").Value = 3837187.5
$ws.Range("L132").Value = 11037.9999
$ws.Range("M132").Value = -3834657.5
$ws.Range("N132").Value = -16097.9999
